# cronograma_curso_provas.xlsx
# "provas gestao de instalacoes e seguranca"
#
# Fill in the exam/unit grades for:
#   row 11 -> "BIM 7D - Gestao de Instalacoes"  (course #9)
#   row 12 -> "BIM 8D - Seguranca"               (course #10)
# and record the class-ranking position (column M) for row 12 and for
# row 13 ("BIM 9D - Lean Construction", course #11), which had been left
# blank even though every other course already has its ranking number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 11: BIM 7D - Gestao de Instalacoes -------------------------------
$ws.Range("F11").Value = 10
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = 10
$ws.Range("I11").Value = 10
$ws.Range("J11").Value = 48
# K11 holds a shared formula (=F11+G11+H11+I11+J11) and recalculates itself.

# --- Row 12: BIM 8D - Seguranca --------------------------------------------
$ws.Range("F12").Value = 10
$ws.Range("G12").Value = 10
$ws.Range("H12").Value = 10
$ws.Range("I12").Value = 10
$ws.Range("J12").Value = 42
# K12 recalculates via its shared formula too.
$ws.Range("L12").Value = "AP"
$ws.Range("M12").Value = 10

# --- Row 13: BIM 9D - Lean Construction ------------------------------------
$ws.Range("M13").Value = 11

# Cursor/selection ends up on M20, matching the saved sheet view.
$ws.Range("M20").Select()
